# Auto-update Data Telemetría - Ejecución Diaria
# Appends the newest daily snapshot (fecha = 45993) to the "historico" sheet
# and replaces the single-row-per-category snapshot on "ultimo_snapshot".

$wb = $excel.ActiveWorkbook

$historico = $wb.Worksheets.Item("historico")
$snapshot  = $wb.Worksheets.Item("ultimo_snapshot")

# Data for the new day (fecha serial 45993 == 2025-12-02)
$rows = @(
    @{ B = "Telemetría";                      C = 5903;  D = 3409; E = 622; F = 191; G = 679; H = 1002; I = 57.75;              J = 10.54;             K = 3.24; L = 11.5;  M = 16.97 },
    @{ B = "GPS (según REGLA)";                C = 5302;  D = 4643; E = 377; F = 100; G = 176; H = 6;    I = 87.56999999999999; J = 7.11;               K = 1.89; L = 3.32;  M = 0.11 },
    @{ B = "GPS (todas con gps_timestamp)";    C = 11199; D = 9377; E = 936; F = 301; G = 585; H = 0;    I = 83.73;              J = 8.359999999999999; K = 2.69; L = 5.22;  M = 0 }
)

$fecha = 45993

# --- Append the 3 new rows to the bottom of "historico" (rows 47-49) ---
$startRow = $historico.Cells.Item($historico.Rows.Count, 1).End(-4162).Row + 1
if ($startRow -lt 2) { $startRow = 2 }

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $historico.Cells.Item($r, 1).Value = $fecha
    $historico.Cells.Item($r, 1).NumberFormat = $historico.Cells.Item($r - 1, 1).NumberFormat

    $historico.Cells.Item($r, 2).Value  = $data.B
    $historico.Cells.Item($r, 3).Value  = $data.C
    $historico.Cells.Item($r, 4).Value  = $data.D
    $historico.Cells.Item($r, 5).Value  = $data.E
    $historico.Cells.Item($r, 6).Value  = $data.F
    $historico.Cells.Item($r, 7).Value  = $data.G
    $historico.Cells.Item($r, 8).Value  = $data.H
    $historico.Cells.Item($r, 9).Value  = $data.I
    $historico.Cells.Item($r, 10).Value = $data.J
    $historico.Cells.Item($r, 11).Value = $data.K
    $historico.Cells.Item($r, 12).Value = $data.L
    $historico.Cells.Item($r, 13).Value = $data.M
}

# --- Overwrite "ultimo_snapshot" rows 2-4 with the same latest data ---
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $data = $rows[$i]

    $snapshot.Cells.Item($r, 1).Value = $fecha

    $snapshot.Cells.Item($r, 2).Value  = $data.B
    $snapshot.Cells.Item($r, 3).Value  = $data.C
    $snapshot.Cells.Item($r, 4).Value  = $data.D
    $snapshot.Cells.Item($r, 5).Value  = $data.E
    $snapshot.Cells.Item($r, 6).Value  = $data.F
    $snapshot.Cells.Item($r, 7).Value  = $data.G
    $snapshot.Cells.Item($r, 8).Value  = $data.H
    $snapshot.Cells.Item($r, 9).Value  = $data.I
    $snapshot.Cells.Item($r, 10).Value = $data.J
    $snapshot.Cells.Item($r, 11).Value = $data.K
    $snapshot.Cells.Item($r, 12).Value = $data.L
    $snapshot.Cells.Item($r, 13).Value = $data.M
}

Write-Output "Done. historico new rows $startRow..$($startRow + $rows.Count - 1); ultimo_snapshot rows 2..4 refreshed."
